# Overview of tags and detections
# Populate No_detections (H) / Plotted (I) columns for several rows, add a
# new Notes (J) value, fix up First/Last detection dates for row 28, and
# add brand new First/Last detection + Exported/No_detections/Plotted data
# for rows 33-35. Also updates the active sheet view/selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 25 ---
$ws.Range("H25").Value = 25
$ws.Range("I25").Value = "Yes"

# --- Row 26 ---
$ws.Range("H26").Value = 4322
$ws.Range("I26").Value = "Yes"
$ws.Range("J26").Value = "Only at Siaes Corner"

# --- Row 27 ---
$ws.Range("H27").Value = 963
$ws.Range("I27").Value = "Yes"

# --- Row 28 ---
$ws.Range("E28").Value = 44867
$ws.Range("F28").Value = 45019
$ws.Range("H28").Value = 3688
$ws.Range("I28").Value = "Yes"

# --- Row 29 ---
$ws.Range("H29").Value = 270
$ws.Range("I29").Value = "Yes"

# --- Row 30 ---
$ws.Range("H30").Value = 24848
$ws.Range("I30").Value = "Yes"

# --- Row 31 ---
$ws.Range("H31").Value = 49
$ws.Range("I31").Value = "Yes"

# --- Row 32 ---
$ws.Range("H32").Value = 2613
$ws.Range("I32").Value = "Yes"

# --- Row 33 (new First_detection/Last_detection/Exported/No_detections/Plotted) ---
$null = $ws.Range("E2").Copy()
$ws.Range("E33").PasteSpecial(-4122)
$ws.Range("E33").Value = 44868
$null = $ws.Range("F2").Copy()
$ws.Range("F33").PasteSpecial(-4122)
$ws.Range("F33").Value = 44990
$ws.Range("G33").Value = "YES"
$ws.Range("H33").Value = 2846
$ws.Range("I33").Value = "Yes"

# --- Row 34 ---
$null = $ws.Range("E2").Copy()
$ws.Range("E34").PasteSpecial(-4122)
$ws.Range("E34").Value = 44868
$null = $ws.Range("F2").Copy()
$ws.Range("F34").PasteSpecial(-4122)
$ws.Range("F34").Value = 44925
$ws.Range("G34").Value = "YES"
$ws.Range("H34").Value = 433
$ws.Range("I34").Value = "Yes"

# --- Row 35 ---
$null = $ws.Range("E2").Copy()
$ws.Range("E35").PasteSpecial(-4122)
$ws.Range("E35").Value = 44869
$null = $ws.Range("F2").Copy()
$ws.Range("F35").PasteSpecial(-4122)
$ws.Range("F35").Value = 45016
$ws.Range("G35").Value = "YES"
$ws.Range("H35").Value = 2402
$ws.Range("I35").Value = "Yes"

$ws.Application.CutCopyMode = $false

# --- Update active selection / scroll position ---
$null = $ws.Range("H28").Select()
